$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.775
$ws.Range("C2").Value = 0.6326530612244898
$ws.Range("D2").Value = 0.6966292134831462
$ws.Range("E2").Value = 49
$ws.Range("B3").Value = 0.660377358490566
$ws.Range("C3").Value = 0.7954545454545454
$ws.Range("D3").Value = 0.7216494845360825
$ws.Range("E3").Value = 44
$ws.Range("B4").Value = 0.7096774193548387
$ws.Range("C4").Value = 0.7096774193548387
$ws.Range("D4").Value = 0.7096774193548387
$ws.Range("E4").Value = 0.7096774193548387
$ws.Range("B5").Value = 0.717688679245283
$ws.Range("C5").Value = 0.7140538033395176
$ws.Range("D5").Value = 0.7091393490096143
$ws.Range("B6").Value = 0.7207699330493
$ws.Range("C6").Value = 0.7096774193548387
$ws.Range("D6").Value = 0.7084667610780838
$ws.Range("B7").Value = 0.7924528301886793
$ws.Range("C7").Value = 0.8571428571428571
$ws.Range("D7").Value = 0.8235294117647058
$ws.Range("E7").Value = 49
$ws.Range("B8").Value = 0.825
$ws.Range("C8").Value = 0.75
$ws.Range("D8").Value = 0.7857142857142856
$ws.Range("E8").Value = 44
$ws.Range("B9").Value = 0.8064516129032258
$ws.Range("C9").Value = 0.8064516129032258
$ws.Range("D9").Value = 0.8064516129032258
$ws.Range("E9").Value = 0.8064516129032258
$ws.Range("B10").Value = 0.8087264150943396
$ws.Range("C10").Value = 0.8035714285714286
$ws.Range("D10").Value = 0.8046218487394957
$ws.Range("B11").Value = 0.8078514911746805
$ws.Range("C11").Value = 0.8064516129032258
$ws.Range("D11").Value = 0.8056383843860125
$ws.Range("B12").Value = 0.8113207547169812
$ws.Range("C12").Value = 0.8775510204081632
$ws.Range("D12").Value = 0.8431372549019608
$ws.Range("E12").Value = 49
$ws.Range("B13").Value = 0.85
$ws.Range("C13").Value = 0.7727272727272727
$ws.Range("D13").Value = 0.8095238095238095
$ws.Range("E13").Value = 44
$ws.Range("B14").Value = 0.8279569892473119
$ws.Range("C14").Value = 0.8279569892473119
$ws.Range("D14").Value = 0.8279569892473119
$ws.Range("E14").Value = 0.8279569892473119
$ws.Range("B15").Value = 0.8306603773584906
$ws.Range("C15").Value = 0.825139146567718
$ws.Range("D15").Value = 0.8263305322128851
$ws.Range("B16").Value = 0.8296206127003449
$ws.Range("C16").Value = 0.8279569892473119
$ws.Range("D16").Value = 0.8272341194542333
$ws.Range("B17").Value = 0.8823529411764706
$ws.Range("C17").Value = 0.9183673469387755
$ws.Range("D17").Value = 0.9
$ws.Range("E17").Value = 49
$ws.Range("B18").Value = 0.9047619047619048
$ws.Range("C18").Value = 0.8636363636363636
$ws.Range("D18").Value = 0.8837209302325582
$ws.Range("E18").Value = 44
$ws.Range("B19").Value = 0.8924731182795699
$ws.Range("C19").Value = 0.8924731182795699
$ws.Range("D19").Value = 0.8924731182795699
$ws.Range("E19").Value = 0.8924731182795699
$ws.Range("B20").Value = 0.8935574229691876
$ws.Range("C20").Value = 0.8910018552875696
$ws.Range("D20").Value = 0.8918604651162791
$ws.Range("B21").Value = 0.8929550314749556
$ws.Range("C21").Value = 0.8924731182795699
$ws.Range("D21").Value = 0.8922980745186295
$ws.Range("B22").Value = 0.8541666666666666
$ws.Range("C22").Value = 0.8367346938775511
$ws.Range("D22").Value = 0.845360824742268
$ws.Range("E22").Value = 49
$ws.Range("B23").Value = 0.8222222222222222
$ws.Range("C23").Value = 0.8409090909090909
$ws.Range("D23").Value = 0.8314606741573033
$ws.Range("E23").Value = 44
$ws.Range("B25").Value = 0.8381944444444445
$ws.Range("C25").Value = 0.8388218923933211
$ws.Range("D25").Value = 0.8384107494497857
$ws.Range("B26").Value = 0.839053166069295
$ws.Range("D26").Value = 0.8387844094117471
